$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6653.8184
$ws.Range("I62").Value = 4721.5
$ws.Range("J62").Value = 7083.222
$ws.Range("K62").Value = 4721.5
$ws.Range("L62").Value = 7083.222
$ws.Range("M62").Value = -4097.5
$ws.Range("N62").Value = -8331.222
$ws.Range("H65").Value = 6653.8184
$ws.Range("I65").Value = 4721.5
$ws.Range("J65").Value = 7083.222
$ws.Range("K65").Value = 23607.5
$ws.Range("L65").Value = 35416.11
$ws.Range("M65").Value = -20487.5
$ws.Range("N65").Value = -41656.11
$ws.Range("H98").Value = 1300.4166
$ws.Range("I98").Value = 1290
$ws.Range("K98").Value = 1290
$ws.Range("M98").Value = 208
$ws.Range("H111").Value = 1396.2354
$ws.Range("I111").Value = 1236.4
$ws.Range("K111").Value = 3709.2
$ws.Range("M111").Value = -642.2000000000003
$ws.Range("H122").Value = 1300.4166
$ws.Range("I122").Value = 1290
$ws.Range("K122").Value = 3870
$ws.Range("M122").Value = -1420
$ws.Range("H132").Value = 2363.9678
$ws.Range("I132").Value = 2335.0386
$ws.Range("K132").Value = 7005.1158
$ws.Range("M132").Value = -4475.1158
$ws.Range("H135").Value = 1158.2106
$ws.Range("I135").Value = 885.61536
$ws.Range("J135").Value = 1748.8334
$ws.Range("K135").Value = 7970.53824
$ws.Range("L135").Value = 15739.5006
$ws.Range("M135").Value = -5435.53824
$ws.Range("N135").Value = -20809.5006
$ws.Range("H138").Value = 5116.586
$ws.Range("J138").Value = 4606.9785
$ws.Range("L138").Value = 13820.9355
$ws.Range("N138").Value = -24100.9355

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2324
$ws.Range("I45").Value = 1648
$ws.Range("K45").Value = 1648
$ws.Range("M45").Value = -1271
$ws.Range("H61").Value = 2398.3333
$ws.Range("I61").Value = 2397.5
$ws.Range("K61").Value = 2397.5
$ws.Range("M61").Value = -2185.5
$ws.Range("H110").Value = 4612.5264
$ws.Range("I110").Value = 4771.6665
$ws.Range("K110").Value = 4771.6665
$ws.Range("M110").Value = -2726.6665
$ws.Range("H136").Value = 2398.3333
$ws.Range("I136").Value = 2397.5
$ws.Range("K136").Value = 7192.5
$ws.Range("M136").Value = -4642.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3295.923
$ws.Range("I107").Value = 1907.5
$ws.Range("K107").Value = 1907.5
$ws.Range("M107").Value = 12.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 33750
$ws.Range("I41").Value = 25000
$ws.Range("J41").Value = 34545.453
$ws.Range("K41").Value = 25000
$ws.Range("L41").Value = 34545.453
$ws.Range("M41").Value = -24572
$ws.Range("N41").Value = -35401.453
$ws.Range("H58").Value = 4368.7144
$ws.Range("I58").Value = 2464.6428
$ws.Range("J58").Value = 8176.857
$ws.Range("K58").Value = 2464.6428
$ws.Range("L58").Value = 8176.857
$ws.Range("M58").Value = -2261.6428
$ws.Range("N58").Value = -8582.857
$ws.Range("H122").Value = 945.5
$ws.Range("I122").Value = 968.25
$ws.Range("K122").Value = 2904.75
$ws.Range("M122").Value = -454.75
$ws.Range("H132").Value = 2636.4348
$ws.Range("I132").Value = 2411.3809
$ws.Range("K132").Value = 7234.1427
$ws.Range("M132").Value = -4704.1427
$ws.Range("H136").Value = 4368.7144
$ws.Range("I136").Value = 2464.6428
$ws.Range("J136").Value = 8176.857
$ws.Range("K136").Value = 7393.928400000001
$ws.Range("L136").Value = 24530.571
$ws.Range("M136").Value = -4843.928400000001
$ws.Range("N136").Value = -29630.571

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 25000404
$ws.Range("I7").Value = 50000250
$ws.Range("K7").Value = 150000750
$ws.Range("M7").Value = -150000638
$ws.Range("H17").Value = 472.5
$ws.Range("I17").Value = 472.5
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 1417.5
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -1248.5
$ws.Range("N17").ClearContents()
$ws.Range("H46").Value = 1250892.8
$ws.Range("I46").Value = 466.25
$ws.Range("J46").Value = 2501319.2
$ws.Range("K46").Value = 1398.75
$ws.Range("L46").Value = 7503957.600000001
$ws.Range("M46").Value = -1307.75
$ws.Range("N46").Value = -7504139.600000001
$ws.Range("H55").Value = 2120
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 2440
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 7320
$ws.Range("M55").Value = -2823
$ws.Range("N55").Value = -7674
$ws.Range("I68").Value = 3032
$ws.Range("J68").Value = 2860.625
$ws.Range("K68").Value = 9096
$ws.Range("L68").Value = 8581.875
$ws.Range("M68").Value = -8285
$ws.Range("N68").Value = -10203.875
$ws.Range("I71").Value = 3032
$ws.Range("J71").Value = 2860.625
$ws.Range("K71").Value = 27288
$ws.Range("L71").Value = 25745.625
$ws.Range("M71").Value = -23232
$ws.Range("N71").Value = -33857.625
$ws.Range("H86").Value = 3250
$ws.Range("J86").Value = 3250
$ws.Range("L86").Value = 9750
$ws.Range("N86").Value = -12122
$ws.Range("H89").Value = 3250
$ws.Range("J89").Value = 3250
$ws.Range("L89").Value = 29250
$ws.Range("N89").Value = -41106
$ws.Range("H92").Value = 415
$ws.Range("I92").Value = 415
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1245
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 3
$ws.Range("N92").ClearContents()
$ws.Range("H102").Value = 3966.6667
$ws.Range("H107").Value = 2476.2222
$ws.Range("J107").Value = 1798
$ws.Range("L107").Value = 5394
$ws.Range("N107").Value = -9234
$ws.Range("H111").Value = 1256.75
$ws.Range("I111").Value = 1256.75
$ws.Range("K111").Value = 3770.25
$ws.Range("M111").Value = -703.25
$ws.Range("I113").Value = 2899.3333
$ws.Range("J113").Value = 3364.8333
$ws.Range("K113").Value = 8697.999899999999
$ws.Range("L113").Value = 10094.4999
$ws.Range("M113").Value = -6527.999899999999
$ws.Range("N113").Value = -14434.4999
$ws.Range("H118").Value = 1679.6
$ws.Range("I118").Value = 1899.5
$ws.Range("J118").Value = 1533
$ws.Range("K118").Value = 5698.5
$ws.Range("L118").Value = 4599
$ws.Range("M118").Value = -4455.5
$ws.Range("N118").Value = -7085
$ws.Range("H129").Value = 2596.2222
$ws.Range("I129").Value = 2845
$ws.Range("J129").Value = 2397.2
$ws.Range("K129").Value = 8535
$ws.Range("L129").Value = 7191.599999999999
$ws.Range("M129").Value = -3535
$ws.Range("N129").Value = -17191.6
$ws.Range("H131").Value = 1391
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 1469.2
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 4407.6
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -14487.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 938.375
$ws.Range("I2").Value = 251.5
$ws.Range("K2").Value = 251.5
$ws.Range("M2").Value = -138.5
$ws.Range("H57").Value = 30998.8
$ws.Range("I57").Value = 22497.5
$ws.Range("J57").Value = 36666.332
$ws.Range("K57").Value = 22497.5
$ws.Range("L57").Value = 36666.332
$ws.Range("M57").Value = -21677.5
$ws.Range("N57").Value = -38306.332
$ws.Range("H80").Value = 5154
$ws.Range("I80").Value = 1502.5
$ws.Range("J80").Value = 7588.3335
$ws.Range("K80").Value = 1502.5
$ws.Range("L80").Value = 7588.3335
$ws.Range("M80").Value = -504.5
$ws.Range("N80").Value = -9584.333500000001
$ws.Range("H83").Value = 5154
$ws.Range("I83").Value = 1502.5
$ws.Range("J83").Value = 7588.3335
$ws.Range("K83").Value = 7512.5
$ws.Range("L83").Value = 37941.6675
$ws.Range("M83").Value = -2520.5
$ws.Range("N83").Value = -47925.6675
$ws.Range("H122").Value = 527832.0600000001
$ws.Range("I122").Value = 70298.92999999999
$ws.Range("K122").Value = 210896.79
$ws.Range("M122").Value = -208446.79

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3728.5
$ws.Range("I46").Value = 2066.6667
$ws.Range("J46").Value = 4974.875
$ws.Range("K46").Value = 2066.6667
$ws.Range("L46").Value = 4974.875
$ws.Range("M46").Value = -1878.6667
$ws.Range("N46").Value = -5350.875
$ws.Range("H68").Value = 4054.5
$ws.Range("J68").Value = 4111
$ws.Range("L68").Value = 4111
$ws.Range("N68").Value = -5609
$ws.Range("H71").Value = 4054.5
$ws.Range("J71").Value = 4111
$ws.Range("L71").Value = 20555
$ws.Range("N71").Value = -28043
$ws.Range("H132").Value = 5169.1577
$ws.Range("I132").Value = 3340.625
$ws.Range("K132").Value = 10021.875
$ws.Range("M132").Value = -7491.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14500
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 14500
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 14500
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -15540
$ws.Range("H81").Value = 3917.5
$ws.Range("J81").Value = 3917.5
$ws.Range("L81").Value = 7835
$ws.Range("N81").Value = -9957
$ws.Range("H84").Value = 3917.5
$ws.Range("J84").Value = 3917.5
$ws.Range("L84").Value = 39175
$ws.Range("N84").Value = -49783
$ws.Range("H107").Value = 2668
$ws.Range("J107").Value = 4000
$ws.Range("L107").Value = 12000
$ws.Range("N107").Value = -15840
$ws.Range("H113").Value = 2334.1667
$ws.Range("I113").Value = 1525.5
$ws.Range("K113").Value = 4576.5
$ws.Range("M113").Value = -2406.5
$ws.Range("H132").Value = 2352.111
$ws.Range("I132").Value = 1718.6154
$ws.Range("K132").Value = 5155.8462
$ws.Range("M132").Value = -2625.8462
